# Auto-generated Excel COM-interop edit script
# Updates cryptos list values (price/volume columns) per commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    # Force the cell to stay a text value: Excel's COM Value setter
    # auto-coerces numeric-looking strings (e.g. "561.02", "1.00") into
    # doubles, which both changes the stored cell type and can introduce
    # floating point drift (e.g. "561.02" -> 561.01999999999998).
    # Marking the cell as Text first preserves the literal string, and
    # resetting the style back to Normal afterwards avoids leaving any
    # visible number-format change behind.
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "64.139.17"
Set-TextValue "E2" "  -0.43%  "
Set-TextValue "D3" "3.148.75"
Set-TextValue "E3" "  +1.98%  "
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "561.02"
Set-TextValue "E5" "  +1.12%  "
Set-TextValue "D6" "140.66"
Set-TextValue "E6" "  +0.71%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "D8" "3.147.03"
Set-TextValue "E8" "  +2.13%  "
Set-TextValue "D9" "0.493"
Set-TextValue "E9" "  -0.01%  "
Set-TextValue "D10" "6.79"
Set-TextValue "E10" "  +3.51%  "
Set-TextValue "E11" "  -0.87%  "
Set-TextValue "E12" "  +0.36%  "
Set-TextValue "D13" "36.19"
Set-TextValue "E13" "  +0.51%  "
Set-TextValue "E14" "  -0.56%  "
Set-TextValue "D15" "3.654.58"
Set-TextValue "E15" "  +2.08%  "
Set-TextValue "D16" "64.144.73"
Set-TextValue "E16" "  -0.56%  "
Set-TextValue "D17" "3.150.92"
Set-TextValue "E17" "  +2.01%  "
Set-TextValue "D19" "510.27"
Set-TextValue "E19" "  +4.09%  "
Set-TextValue "D20" "6.78"
Set-TextValue "E20" "  +1.43%  "
Set-TextValue "E21" "  +1.79%  "
Set-TextValue "D22" "0.713"
Set-TextValue "E22" "  +3.24%  "
Set-TextValue "E23" "  +2.26%  "
Set-TextValue "D24" "12.70"
Set-TextValue "E24" "  +2.03%  "
Set-TextValue "D25" "78.69"
Set-TextValue "E25" "  -0.08%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "D27" "8.67"
Set-TextValue "E27" "  +7.90%  "
Set-TextValue "E28" "  +2.83%  "
Set-TextValue "E29" "  +0.30%  "
Set-TextValue "E30" "  -0.14%  "
Set-TextValue "D31" "26.59"
Set-TextValue "E31" "  +1.84%  "
Set-TextValue "E32" "  -1.53%  "
Set-TextValue "E33" "  -0.63%  "
Set-TextValue "D34" "554.24"
Set-TextValue "E34" "  -4.62%  "
Set-TextValue "D35" "6.05"
Set-TextValue "E35" "  +0.13%  "
Set-TextValue "D36" "53.86"
Set-TextValue "E36" "  +2.19%  "
Set-TextValue "D37" "5.31"
Set-TextValue "E37" "  -2.82%  "
Set-TextValue "E38" "  +4.23%  "
Set-TextValue "D39" "3.150.22"
Set-TextValue "E39" "  +6.25%  "
Set-TextValue "E40" "  +2.27%  "
Set-TextValue "E41" "  +1.61%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D42" "8.23"
Set-TextValue "E42" "  -1.02%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.71"
Set-TextValue "E43" "  -8.29%  "
Set-TextValue "E44" "  +7.09%  "
Set-TextValue "D45" "2.16"
Set-TextValue "E45" "  +1.85%  "
Set-TextValue "D46" "1.00"
Set-TextValue "E46" "  +0.00%  "
Set-TextValue "D47" "122.34"
Set-TextValue "E47" "  +1.78%  "
Set-TextValue "D48" "24.83"
Set-TextValue "E48" "  -2.23%  "
Set-TextValue "E49" "  -0.75%  "
Set-TextValue "D50" "0.0₃0512"
Set-TextValue "E50" "  -5.42%  "
Set-TextValue "E51" "  -0.30%  "
